$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - force text format to preserve exact formatting
# (avoids Excel auto-converting numeric-looking strings and dropping trailing zeros)
$priceUpdates = [ordered]@{
    'D2' = '61.798.49'
    'D3' = '2.491.62'
    'D5' = '554.75'
    'D6' = '147.32'
    'D9' = '2.491.67'
    'D10' = '0.108'
    'D11' = '5.45'
    'D14' = '26.44'
    'D15' = '2.940.56'
    'D17' = '61.719.25'
    'D18' = '2.489.89'
    'D20' = '7.03'
    'D21' = '4.23'
    'D22' = '324.09'
    'D24' = '1.81'
    'D25' = '64.26'
    'D26' = '0.0₃0999'
    'D28' = '2.612.07'
    'D29' = '1.00'
    'D34' = '1.92'
    'D35' = '1.61'
    'D36' = '6.00'
    'D37' = '4.94'
    'D38' = '0.999'
    'D39' = '0.385'
    'D40' = '18.63'
    'D41' = '148.45'
    'D43' = '0.999'
    'D44' = '40.42'
    'D45' = '2.35'
    'D46' = '149.38'
    'D48' = '21.13'
    'D49' = '0.0540'
}
foreach ($addr in $priceUpdates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $priceUpdates[$addr]
}

# Volume(1h) (column E) updates
$volumeUpdates = [ordered]@{
    'E2' = '  -3.01%  '
    'E3' = '  -5.26%  '
    'E4' = '  +0.04%  '
    'E5' = '  -3.84%  '
    'E6' = '  -5.01%  '
    'E7' = '  -0.01%  '
    'E8' = '  -3.11%  '
    'E9' = '  -5.17%  '
    'E10' = '  -7.84%  '
    'E11' = '  -6.11%  '
    'E12' = '  -1.39%  '
    'E13' = '  -5.13%  '
    'E14' = '  -6.81%  '
    'E15' = '  -5.33%  '
    'E16' = '  -7.25%  '
    'E17' = '  -3.06%  '
    'E18' = '  -5.86%  '
    'E19' = '  -7.32%  '
    'E20' = '  -7.88%  '
    'E21' = '  -6.45%  '
    'E22' = '  -5.93%  '
    'E24' = '  -4.06%  '
    'E25' = '  -5.41%  '
    'E26' = '  -7.71%  '
    'E27' = '  -3.60%  '
    'E28' = '  -5.08%  '
    'E29' = '  +0.05%  '
    'E30' = '  -10.92%  '
    'E31' = '  -8.85%  '
    'E32' = '  -4.21%  '
    'E33' = '  -5.31%  '
    'E34' = '  -7.14%  '
    'E35' = '  -8.24%  '
    'E36' = '  -9.57%  '
    'E37' = '  -8.41%  '
    'E38' = '  +0.02%  '
    'E39' = '  -4.15%  '
    'E40' = '  -5.53%  '
    'E41' = '  -1.08%  '
    'E42' = '  -7.79%  '
    'E43' = '  +0.03%  '
    'E44' = '  -3.15%  '
    'E45' = '  -6.96%  '
    'E46' = '  -6.17%  '
    'E47' = '  -6.70%  '
    'E48' = '  -14.74%  '
    'E49' = '  -7.89%  '
    'E50' = '  -4.83%  '
    'E51' = '  -4.68%  '
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}

Write-Host "Updated $($priceUpdates.Count) price cells and $($volumeUpdates.Count) volume cells."
